$d = $word.ActiveDocument

# --- 1) "Informasi lain Nomor : ${sumber_npi} tanggal ${tgl_li}"
#        -> "Informasi lain Nomor : ${sumber_informasi} "
$rng = $d.Content
$rng.Find.Execute('Informasi lain Nomor : ${sumber_npi} tanggal ${tgl_li}', $false, $false, $false, $false, $false, $true, 1, $false, 'Informasi lain Nomor : ${sumber_informasi} ', 2) | Out-Null

# --- 2) "Kategori Penindakan : " -> "Kategori Penindakan : ${kategori_npi}"
$rng = $d.Content
$rng.Find.Execute('Kategori Penindakan : ', $false, $false, $false, $false, $false, $true, 1, $false, 'Kategori Penindakan : ${kategori_npi}', 2) | Out-Null

# --- 3) remove the stray "Uni" before " ${unit_penerbit_npi}"
$rng = $d.Content
$rng.Find.Execute('informasi dari Uni ${unit_penerbit_npi}', $false, $false, $false, $false, $false, $true, 1, $false, 'informasi dari ${unit_penerbit_npi}', 2) | Out-Null

# --- 4) add a new empty (justified) paragraph between "...berikut:" and the table
$rng = $d.Content
$rng.Find.Execute('belum* dapat dilakukan Penindakan lebih lanjut dengan alasan sebagai berikut:', $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertBefore("`r")

# --- 5) widen the table / grid / cells from 9350 dxa to 9493 dxa, and justify the
#        first data cell's paragraph
$t = $d.Tables(1)
$t.PreferredWidthType = 3
$t.PreferredWidth = 474.65
$t.Columns(1).Width = 474.65
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $t.Cell($r, 1).Width = 474.65
}
$t.Cell(1, 1).Range.Paragraphs(1).Alignment = 3

# --- 6) add 3 new empty paragraphs right after "...terima kasih."
for ($i = 0; $i -lt 3; $i++) {
    $rng = $d.Content
    $rng.Find.Execute('kasih.', $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.Collapse(0)
    $rng.InsertBefore("`r")
}
